# ticker.history period changed to 3d, keyword for record skipping changed
# from 'PLN' to 'none'.

$wb = $excel.ActiveWorkbook

$assets = $wb.Worksheets.Item("assets")

# --- "assets" sheet data updates -----------------------------------------

# BTC-USD price refreshed after re-pulling a shorter (3 day) history window.
$assets.Range("D5").Value = 108106.3671875

# Rows that used to be flagged with the currency ticker "PLN" (a sentinel
# meaning "this isn't a real market ticker, skip it when pulling history")
# now use the literal keyword "none" instead.
$assets.Range("A6").Value = "none"
$assets.Range("A7").Value = "none"

# --- Active sheet/tab switched to "assets" --------------------------------
$assets.Select()
